$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 301, pushing existing rows 301-371 down to 303-373
$ws.Rows("301:302").Insert()

# Populate new row 301
$ws.Cells.Item(301, 1).Value = 9
$ws.Cells.Item(301, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(301, 3).Value = "Metropolitana"
$ws.Cells.Item(301, 4).Value = 44943
$ws.Cells.Item(301, 5).Value = 13
$ws.Cells.Item(301, 6).Value = 100112017
$ws.Cells.Item(301, 7).Value = "Apio"
$ws.Cells.Item(301, 8).Value = "Americana (o)"
$ws.Cells.Item(301, 9).Value = "Primera"
$ws.Cells.Item(301, 10).Value = 70
$ws.Cells.Item(301, 11).Value = 8000
$ws.Cells.Item(301, 12).Value = 9000
$ws.Cells.Item(301, 13).Value = 8500
$ws.Cells.Item(301, 14).Value = "`$/docena de matas"
$ws.Cells.Item(301, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(301, 16).Value = 1417
$ws.Cells.Item(301, 17).Value = 6
$ws.Cells.Item(301, 18).Value = "Hortaliza"

# Populate new row 302
$ws.Cells.Item(302, 1).Value = 9
$ws.Cells.Item(302, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(302, 3).Value = "Metropolitana"
$ws.Cells.Item(302, 4).Value = 44943
$ws.Cells.Item(302, 5).Value = 13
$ws.Cells.Item(302, 6).Value = 100112017
$ws.Cells.Item(302, 7).Value = "Apio"
$ws.Cells.Item(302, 8).Value = "Americana (o)"
$ws.Cells.Item(302, 9).Value = "Segunda"
$ws.Cells.Item(302, 10).Value = 43
$ws.Cells.Item(302, 11).Value = 7000
$ws.Cells.Item(302, 12).Value = 7000
$ws.Cells.Item(302, 13).Value = 7000
$ws.Cells.Item(302, 14).Value = "`$/docena de matas"
$ws.Cells.Item(302, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(302, 16).Value = 1167
$ws.Cells.Item(302, 17).Value = 6
$ws.Cells.Item(302, 18).Value = "Hortaliza"
